$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 and J1, matching the style of the existing header cells (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF) for rows 2 through 62
$data = @(
    @(2; 9; 9),
    @(3; 8; 9),
    @(4; 9; 9),
    @(5; 9; 9),
    @(6; 9; 9),
    @(7; 9; 9),
    @(8; 9; 9),
    @(9; 9; 9),
    @(10; 8; 8),
    @(11; 9; 9),
    @(12; 9; 9),
    @(13; 9; 9),
    @(14; 9; 9),
    @(15; 9; 9),
    @(16; 10; 10),
    @(17; 9; 9),
    @(18; 9; 9),
    @(19; 9; 9),
    @(20; 9; 9),
    @(21; 8; 9),
    @(22; 9; 9),
    @(23; 7; 8),
    @(24; 8; 9),
    @(25; 7; 7),
    @(26; 8; 9),
    @(27; 9; 9),
    @(28; 8; 9),
    @(29; 9; 9),
    @(30; 9; 9),
    @(31; 9; 9),
    @(32; 8; 9),
    @(33; 9; 9),
    @(34; 9; 9),
    @(35; 8; 8),
    @(36; 9; 9),
    @(37; 9; 9),
    @(38; 8; 8),
    @(39; 8; 8),
    @(40; 7; 7),
    @(41; 7; 7),
    @(42; 7; 8),
    @(43; 9; 9),
    @(44; 7; 8),
    @(45; 6; 6),
    @(46; 8; 8),
    @(47; 7; 7),
    @(48; 8; 8),
    @(49; 8; 8),
    @(50; 8; 9),
    @(51; 9; 9),
    @(52; 9; 9),
    @(53; 8; 9),
    @(54; 9; 9),
    @(55; 6; 6),
    @(56; 7; 7),
    @(57; 6; 6),
    @(58; 4; 5),
    @(59; 3; 4),
    @(60; 7; 7),
    @(61; 4; 4),
    @(62; 6; 6)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
